$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# The department column (C) previously held the single shared string
# "FACULTY OF BUSINESS & TECHNOLOGY" for every course row. Replace it with
# the per-course department name: the two Information Technology courses
# (rows 2-3) get "Information Technology" and the Civil Construction course
# (row 4) gets "Building and Construction".
$ws.Range("C2").Value = "Information Technology"
$ws.Range("C3").Value = "Information Technology"
$ws.Range("C4").Value = "Building and Construction"

# Leave the selection on the last-edited cell, matching the editor's
# recorded cursor position after making the change.
[void]$ws.Range("C4").Select()
